$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "14.77%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.082"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.18%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07809"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.63%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.276"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.48%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.097"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.19%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.02%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9282"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.06%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09763"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.69%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1818"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.45%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08689"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.21%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03415"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.33%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09926"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.01%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001472"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.84%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005721"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.06%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.485"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.52%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.11%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3432"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.02%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.01%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.551"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "11.31%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004492"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.28%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.10%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-20.32%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01753"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.43%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04696"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.47%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007859"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.30%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.45%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008574"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-12.83%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002301"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.71%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009184"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.80%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006131"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.45%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.11%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.056"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "52.80%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.78%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.11%"
